$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.440.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3797"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07888"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9682"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.54%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.832.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.870"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.054"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06633"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001025"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.435.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.326"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.298"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.027.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.058"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9395"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09292"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.584"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.241"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05929"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02177"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.055"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.142"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5758"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1825"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.966"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.266"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5442"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.865"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06608"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.041"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.44%  "
